$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" (columns A:N) ---
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row at position 220 (shifts existing rows 220-259 down to 221-260)
$ws1.Rows.Item(220).Insert()

# Populate the newly inserted row with the new client record "DISALME CIA. LTDA."
# under advisor "LOZANO MOLINA TITO" (same advisor as the surrounding rows), all
# numeric columns at zero.
$ws1.Range("A220").Value = "LOZANO MOLINA TITO"
$ws1.Range("B220").Value = "DISALME CIA. LTDA."
$ws1.Range("C220").Value = 0
$ws1.Range("D220").Value = 0
$ws1.Range("E220").Value = 0
$ws1.Range("F220").Value = 0
$ws1.Range("G220").Value = 0
$ws1.Range("H220").Value = 0
$ws1.Range("I220").Value = 0
$ws1.Range("J220").Value = 0
$ws1.Range("K220").Value = 0
$ws1.Range("L220").Value = 0
$ws1.Range("M220").Value = 0
$ws1.Range("N220").Value = 0

# Update the "N de <count>" summary row (now row 260) to reflect the new total
# row count of 258 (was 257).
$ws1.Range("C260").Value = "2 de 258"
$ws1.Range("D260").Value = "24 de 258"
$ws1.Range("E260").Value = "5 de 258"
$ws1.Range("F260").Value = "1 de 258"
$ws1.Range("G260").Value = "4 de 258"
$ws1.Range("H260").Value = "3 de 258"
$ws1.Range("I260").Value = "2 de 258"
$ws1.Range("J260").Value = "5 de 258"
$ws1.Range("K260").Value = "21 de 258"
$ws1.Range("L260").Value = "57 de 258"
$ws1.Range("M260").Value = "3 de 258"
$ws1.Range("N260").Value = "5 de 258"

# --- Sheet 2: "VENTA MENSUAL" (columns A:F) ---
$ws2 = $wb.Worksheets.Item(2)

# Same insert, mirrored onto the second sheet.
$ws2.Rows.Item(220).Insert()

$ws2.Range("A220").Value = "LOZANO MOLINA TITO"
$ws2.Range("B220").Value = "DISALME CIA. LTDA."
$ws2.Range("C220").Value = 0
$ws2.Range("D220").Value = 0
$ws2.Range("E220").Value = 0
$ws2.Range("F220").Value = 0
